# Update chart code and README
# Renames the mark-sheet subject headers and corrects a batch of mark
# values that feed the chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1: subject columns) ---
$ws.Range("C1").Value = "ADMS"
$ws.Range("D1").Value = "AOS"
$ws.Range("E1").Value = "A&CD"
$ws.Range("F1").Value = "C&NS"

# --- Corrected mark values ---
$ws.Range("F5").Value = 14
$ws.Range("D6").Value = 12
$ws.Range("F8").Value = 16
$ws.Range("D9").Value = 16
$ws.Range("D14").Value = 13
$ws.Range("C15").Value = 12
$ws.Range("C18").Value = 11
$ws.Range("D18").Value = 18
$ws.Range("C19").Value = 90
$ws.Range("E19").Value = 89
$ws.Range("D21").Value = 17
$ws.Range("C22").Value = 14
$ws.Range("D27").Value = 17
$ws.Range("F29").Value = 12
$ws.Range("F30").Value = 16
$ws.Range("C31").Value = 98
$ws.Range("F31").Value = 94
$ws.Range("E33").Value = 18
$ws.Range("D34").Value = 13
$ws.Range("D35").Value = 18
$ws.Range("E35").Value = 17
$ws.Range("F35").Value = 18
$ws.Range("D36").Value = 15
$ws.Range("E36").Value = 17
$ws.Range("D38").Value = 13
$ws.Range("C40").Value = 91
$ws.Range("E40").Value = 86
$ws.Range("F40").Value = 96
$ws.Range("E41").Value = 17

# --- Selection moved to F1 to mirror the saved view state ---
$ws.Range("F1").Select()
